# ============================================================================
# Adds a "Player Info" sheet (player bio) and an "ODI Batting Extra" sheet
# (per-match batting-position / boundary-count detail), and converts the
# MATCH_CARD_LINK url column on the existing "ODI Batting" / "ODI Bowling"
# sheets into a plain MATCH_CODE column.
# ============================================================================

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Helper: apply the workbook's existing bold/border/center-top header look
# to a range (mirrors the style already used on row 1 of the other sheets).
# ----------------------------------------------------------------------
function Format-HeaderRange($range) {
    $range.Font.Bold = $true
    $range.HorizontalAlignment = -4108   # xlCenter
    $range.VerticalAlignment = -4160     # xlTop
    $range.Borders.LineStyle = 1         # xlContinuous
}

# ----------------------------------------------------------------------
# 1. Rewrite "ODI Batting": MATCH_CARD_LINK -> MATCH_CODE (text code only)
# ----------------------------------------------------------------------
$wsBatting = $wb.Worksheets.Item("ODI Batting")
$wsBatting.Range("D1").Value = "MATCH_CODE"

$lastRowBatting = $wsBatting.Cells.Item($wsBatting.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRowBatting; $r++) {
    $cell = $wsBatting.Cells.Item($r, 4)
    $url = $cell.Value2
    $code = $url -replace ".*MatchCode=", ""
    $cell.NumberFormat = "@"
    $cell.Value = $code
}

# ----------------------------------------------------------------------
# 2. Rewrite "ODI Bowling": MATCH_CARD_LINK -> MATCH_CODE (text code only)
# ----------------------------------------------------------------------
$wsBowling = $wb.Worksheets.Item("ODI Bowling")
$wsBowling.Range("B1").Value = "MATCH_CODE"

$lastRowBowling = $wsBowling.Cells.Item($wsBowling.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRowBowling; $r++) {
    $cell = $wsBowling.Cells.Item($r, 2)
    $url = $cell.Value2
    $code = $url -replace ".*MatchCode=", ""
    $cell.NumberFormat = "@"
    $cell.Value = $code
}

# ----------------------------------------------------------------------
# 3. Insert new "Player Info" sheet as the first tab
# ----------------------------------------------------------------------
$wsPlayer = $wb.Worksheets.Add($wb.Worksheets.Item(1))
$wsPlayer.Name = "Player Info"

$wsPlayer.Range("A1").Value = "ID"
$wsPlayer.Range("B1").Value = "NAME"
$wsPlayer.Range("C1").Value = "BATTING_HAND"
$wsPlayer.Range("D1").Value = "BOWL_STYLE"
Format-HeaderRange($wsPlayer.Range("A1:D1"))

$wsPlayer.Cells.Item(2, 1).NumberFormat = "@"
$wsPlayer.Cells.Item(2, 1).Value = "3269"
$wsPlayer.Cells.Item(2, 2).Value = "Xavier Melbourne Marshall"
$wsPlayer.Cells.Item(2, 3).Value = "Right Handed"
$wsPlayer.Cells.Item(2, 4).Value = "Right Arm Off Break"

# ----------------------------------------------------------------------
# 4. Append new "ODI Batting Extra" sheet as the last tab
# ----------------------------------------------------------------------
$wsBowling = $wb.Worksheets.Item("ODI Bowling")
$wsExtra = $wb.Worksheets.Add($null, $wsBowling)
$wsExtra.Name = "ODI Batting Extra"

$wsExtra.Range("A1").Value = "MATCH_CODE"
$wsExtra.Range("B1").Value = "BATTING_POSITION"
$wsExtra.Range("C1").Value = "NUM_4"
$wsExtra.Range("D1").Value = "NUM_6"
$wsExtra.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$wsExtra.Range("F1").Value = "MAN_OF_MATCH"
Format-HeaderRange($wsExtra.Range("A1:F1"))

# MATCH_CODE, BATTING_POSITION, NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH
$extraData = @(
    @("2867","5","1","0","2.21%","NO"),
    @("2894","5","4","0","8.84%","NO"),
    @("2897","5","0","0","0.48%","NO"),
    @("2898","5","0","0","","NO"),
    @("2908","","","","","NO"),
    @("2909","","","","","NO"),
    @("2910","2","0","0","","NO"),
    @("4282","1","0","1","10.66%","NO"),
    @("4368","","","","","NO"),
    @("4369","2","2","1","16.50%","NO"),
    @("4370","","","","","NO"),
    @("4371","","","","","NO"),
    @("4380","1","0","0","1.29%","NO"),
    @("4381","1","2","1","12.50%","NO"),
    @("4383","","","","",""),
    @("4384","","","","",""),
    @("4404","","","","",""),
    @("4407","","","","",""),
    @("4411","","","","",""),
    @("4412","","","","","")
)

$r = 2
foreach ($row in $extraData) {
    $wsExtra.Cells.Item($r, 1).NumberFormat = "@"
    $wsExtra.Cells.Item($r, 1).Value = $row[0]

    if ($row[1] -ne "") {
        $wsExtra.Cells.Item($r, 2).Value = [int]$row[1]
    }
    if ($row[2] -ne "") {
        $wsExtra.Cells.Item($r, 3).NumberFormat = "@"
        $wsExtra.Cells.Item($r, 3).Value = $row[2]
    }
    if ($row[3] -ne "") {
        $wsExtra.Cells.Item($r, 4).NumberFormat = "@"
        $wsExtra.Cells.Item($r, 4).Value = $row[3]
    }
    if ($row[4] -ne "") {
        $wsExtra.Cells.Item($r, 5).NumberFormat = "@"
        $wsExtra.Cells.Item($r, 5).Value = $row[4]
    }
    if ($row[5] -ne "") {
        $wsExtra.Cells.Item($r, 6).Value = $row[5]
    }
    $r = $r + 1
}

# ----------------------------------------------------------------------
# 5. Restore the active tab to the first sheet (matches original workbook
#    view state: activeTab="0"/firstSheet="0").
# ----------------------------------------------------------------------
$wb.Worksheets.Item(1).Activate()

Write-Output "Edit complete"
